$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("User_Initiated_Messages")
$ws4 = $wb.Worksheets.Item("Follow_Up_Messages")

# Fix literal "\n" escape sequences in the bot response message text
# (double backslash-n -> plain spaces / semicolons)
$ws4.Range("C2").Value = 'Great, there''s lots of benefits by not gambling.  Some benefits are getting out of debt, improving relationships, feeling healthier and less stressed, feeling better about yourself.  What do you see as your benefits by stopping gambling?'
$ws4.Range("C5").Value = 'Great, here are some quick tips to start this journey:  Avoid tempting environments; Limit access to finances; Find healthier activities to do'
$ws4.Range("C7").Value = 'Try to limit how many days you go gambling, or set a maximum amount you can spend.  Take some time now to write down what you want to achieve'

# Fix the regex key strings: collapse the erroneous double backslash
# ("\\b") down to a single backslash ("\b") word-boundary escape
$ws4.Range("B2").Value = '^(?=.*\byes\b.*).*$'
$ws4.Range("B5").Value = '^(?=.*\byes\b.*).*$'
$ws4.Range("B7").Value = '^(?=.*\byes\b.*).*$'
$ws4.Range("B3").Value = '^(?=.*\bno\b.*).*$'
$ws4.Range("B6").Value = '^(?=.*\bno\b.*).*$'
$ws4.Range("B8").Value = '^(?=.*\bno\b.*).*$'
$ws3.Range("B2").Value = '^(?=.*\bgoals?\b.*).*$'

Write-Host "Done applying fixes"
